$d = $word.ActiveDocument

# --- Hunk 1: insert new paragraph about computational complexity (with two OMath expressions) ---
# The target paragraph is the one describing the pivot selection ("Podstawowym problemem...").
$targetIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx++
    if ($p.Range.Text -like "*Podstawowym problemem jest*") {
        $targetIdx = $idx
    }
}
if ($targetIdx -eq -1) {
    throw "Hunk1: target paragraph not found"
}
$para = $d.Paragraphs($targetIdx)
$r = $para.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$newPara = $d.Paragraphs($targetIdx + 1)
$newPara.Range.Text = "PLACEHOLDER_FOR_NEW_PARAGRAPH"
$mrng = $newPara.Range
$hunk1Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Z</w:t></w:r><w:r><w:t>łożoność obliczeniowa</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>O</m:t></m:r><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:dPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>n</m:t></m:r><m:func><m:funcPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:funcPr><m:fName><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>log</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:sub></m:sSub></m:fName><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>n</m:t></m:r></m:e></m:func></m:e></m:d></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>dla dobrze wybranego „</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pivota</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">i </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>O</m:t></m:r><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:dPr><m:e><m:sSup><m:sSupPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSupPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>n</m:t></m:r></m:e><m:sup><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> dla wybranego źle.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$mrng.InsertXML($hunk1Xml)
Write-Host "Hunk1 done"

# --- Hunk 2: merge "pivot_index=" and "1" math runs into a single run (pivot_index=1) ---
# Target paragraph contains "Ponownie trafiamy na wartosc mniejsza" description.
$targetIdx2 = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx++
    if ($p.Range.Text -like "*Ponownie trafiamy*") {
        $targetIdx2 = $idx
    }
}
if ($targetIdx2 -eq -1) {
    throw "Hunk2: target paragraph not found"
}
$target2 = $d.Paragraphs($targetIdx2)
$hunk2Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16du="http://schemas.microsoft.com/office/word/2023/wordml/word16du" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16sdtfl="http://schemas.microsoft.com/office/word/2024/wordml/sdtformatlock" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh w16sdtfl w16du wp14"><w:body><w:p w14:paraId="6E4E6606" w14:textId="6294D05A" w:rsidR="00A56D31" w:rsidRDefault="00A56D31" w:rsidP="00A56D31"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">W pierwszym kroku algorytmu, analizowana będzie cała tablica. Na początku losujemy indeks </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pivota</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> – przyjmijmy, że wylosowała się wartość </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:oMath><w:r><w:t>, czyli</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>pivot_value=</m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:shd w:val="clear" w:color="auto" w:fill="E59EDC" w:themeFill="accent5" w:themeFillTint="66"/></w:rPr><m:t>5</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">. Aby uniknąć przypadkowego przestawienia elementu </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>pivota</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">, przeniesiemy go tymczasowo na koniec tablicy. Następnie będziemy sprawdzać kolejne wartości od strony lewej i porównywać je z wartością </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>pivota</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">. Jeżeli sprawdzana wartość będzie mniejsza, zamienimy ją z wartością stojącą na pozycji </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>pivot_index</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> (czyli na indeksie odpowiadającym ilości wartości mniejszych), którą zwiększymy. W naszym przypadku najpierw natrafiamy na liczbę </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:shd w:val="clear" w:color="auto" w:fill="95DCF7" w:themeFill="accent4" w:themeFillTint="66"/></w:rPr><m:t>7</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">, która jest od pivota większa, więc przechodzimy od razu do kolejnej wartości. Liczba </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:shd w:val="clear" w:color="auto" w:fill="95DCF7" w:themeFill="accent4" w:themeFillTint="66"/></w:rPr><m:t>3</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> jest mniejsza od </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:shd w:val="clear" w:color="auto" w:fill="E59EDC" w:themeFill="accent5" w:themeFillTint="66"/></w:rPr><m:t>5</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">, dlatego zamieniamy ją z elementem stojącym na pozycji </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>pivot_index=0</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> oraz zwiększamy licznik o </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>1</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">. Ponownie trafiamy na wartość mniejszą, </w:t></w:r><w:r w:rsidRPr="00A56D31"><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:shd w:val="clear" w:color="auto" w:fill="95DCF7" w:themeFill="accent4" w:themeFillTint="66"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">, którą zamieniamy z elementem stojącym na pozycji </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>pivot_index=1</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> i zwiększamy licznik. Ostatnią sprawdzaną wartością jest </w:t></w:r><w:r w:rsidRPr="00A56D31"><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:shd w:val="clear" w:color="auto" w:fill="95DCF7" w:themeFill="accent4" w:themeFillTint="66"/></w:rPr><w:t>9</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">, która jest większa od wartości </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>pivota</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">, dlatego nie wykonujemy żadnych operacji i wychodzimy z pętli. Kluczowym krokiem jest przestawienie teraz z powrotem wartości </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>pivota</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> we właściwe miejsce – na pozycję </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>pivot_index</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">. Dwie powstałe części – </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>[3, 3]</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> oraz </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>[9, 7]</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>, które są od niego odpowiednio mniejsze oraz niemniejsze, posortujemy teraz wykorzystując algorytm sortowania szybkiego (a więc rekurencyjnie).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target2.Range.InsertXML($hunk2Xml)
Write-Host "Hunk2 done"

# --- Hunk 3: merge "Tabela 3 Krok 3. ..." runs into a single run ---
$find = $d.Content.Find
$found3 = $find.Execute("Tabela 3 Krok 3.", $false, $false, $false, $false, $false, $true, 1, $false, "Tabela 3 Krok 3.", 2)
Write-Host "Hunk3 found=$found3"
if (-not $found3) {
    throw "Hunk3: target text not found"
}
Write-Host "All done"
